# S sensitivity added to the function
# Refresh the computed sensitivity values on the active sheet
# ("CO2 Change by Activities") with the newly recalculated results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = 0.00002232981138561296
    "D2"  = 0.004400939973805862
    "F2"  = 0.00002232981138561296
    "G2"  = 0.0001116490569281758

    "C3"  = 0.000000003276756976866579
    "D3"  = 0.00000001310702790746632
    "E3"  = 0.000000004006398324207794
    "F3"  = 0.000000006073774150294753
    "G3"  = 0.0000000163837849953552

    "C4"  = 0.00000000001054685679069411
    "D4"  = 0.00000000004218742716277646
    "E4"  = 0.00000004823850474622304
    "F4"  = 0.000000007800784906919489
    "G4"  = 0.0000000000527342822187471

    "C5"  = 0.07578661345678483
    "D5"  = 0.0001272919309522891
    "F5"  = 0.00003182298273807227
    "G5"  = 0.0001591149136901393

    "C6"  = 0.00000001478532157817014
    "D6"  = 0.000002914010849508486
    "F6"  = 0.00000001478532157817014
    "G6"  = 0.00000007392660783533955

    "C7"  = 0.00000001930840100305975
    "D7"  = 0.000000077233604012239
    "E7"  = 0.0000883114652765471
    "F7"  = 0.00001428109840162506
    "G7"  = 0.00000009654200505693211

    "B8"  = 0.00330747736916237
    "C8"  = 0.0004874692202960773
    "D8"  = 0.08780026868078039
    "F8"  = 0.0004874692202960773
    "G8"  = 0.002437346101480387

    "C9"  = 0.00001342734235232346
    "D9"  = 0.01156986263652016
    "F9"  = 0.00001342734235232346
    "G9"  = 0.00006713671176150626

    "C10" = 0.00003994091957792989
    "D10" = 0.003533260433755459
    "E10" = 0.00001500819767841222
    "F10" = 0.00005412836907225937
    "G10" = 0.0001997045978896494

    "C11" = 0.000007279715729779923
    "D11" = 0.001560155710080835
    "F11" = 0.0006627228178786027
    "G11" = 0.01782423959866719
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
